# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The column "K" (column G, header "K") holds new values for each row (2-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @(3, 1, 0, 7, 5, 4, 4, 2, 3, 1, 2, 2, 2, 1, 0, 2, 0, 0, 1, 1, 0, 0, 4, 0, 0, 4, 4, 1, 3, 4, 3, 3, 4, 5, 2, 5, 5, 5, 4, 5, 3, 2, 3)

$row = 2
foreach ($val in $newKValues) {
    $ws.Range("G$row").Value = $val
    $row++
}
